$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new title row at the top; existing rows shift down by one.
$ws.Rows(1).Insert()

# Title cell: "Our cool Report", merged across A1:C1, big blue centered text.
$ws.Range("A1").Value = "Our cool Report"
$ws.Range("A1:C1").Merge()
$ws.Range("A1").Font.Size = 24
$ws.Range("A1").Font.Color = 16711680
$ws.Range("A1").HorizontalAlignment = -4108

# Header row (was row 1, now row 2): Id / FirstName / LastName -> bold + centered
$ws.Range("A2:C2").Font.Bold = $true
$ws.Range("A2:C2").HorizontalAlignment = -4108

# Data rows: center the Id column (column A) for every data row
$ws.Range("A3:A5").HorizontalAlignment = -4108

# Column C is widened to fit the longer text
$ws.Columns("C").ColumnWidth = 19.166666666666668

Write-Host "done"
